{"js": "const replacements = [\n  [\"2024-08-25 Sunday\", \"2024-08-26 Monday\"],\n  [\"37\u00f76=6, 1\", \"31\u00f76=5, 1\"],\n  [\"67\u00f79=7, 4\", \"57\u00f72=28, 1\"],\n  [\"10\u00f75=2, 0\", \"43\u00f75=8, 3\"],\n  [\"60\u00f76=10, 0\", \"93\u00f77=13, 2\"],\n  [\"66\u00f72=33, 0\", \"14\u00f76=2, 2\"],\n  [\"30\u00f76=5, 0\", \"97\u00f73=32, 1\"],\n  [\"88\u00f76=14, 4\", \"60\u00f77=8, 4\"],\n  [\"69\u00f77=9, 6\", \"44\u00f72=22, 0\"],\n  [\"15\u00f79=1, 6\", \"39\u00f76=6, 3\"],\n  [\"75\u00f74=18, 3\", \"42\u00f76=7, 0\"],\n  [\"98\u00f73=32, 2\", \"32\u00f75=6, 2\"],\n  [\"24\u00f73=8, 0\", \"95\u00f74=23, 3\"],\n  [\"79\u00f73=26, 1\", \"33\u00f78=4, 1\"],\n  [\"81\u00f79=9, 0\", \"13\u00f72=6, 1\"],\n  [\"43\u00f72=21, 1\", \"12\u00f76=2, 0\"],\n  [\"99\u00f72=49, 1\", \"47\u00f72=23, 1\"],\n  [\"51\u00f74=12, 3\", \"71\u00f78=8, 7\"],\n  [\"77\u00f72=38, 1\", \"65\u00f79=7, 2\"],\n  [\"43\u00f77=6, 1\", \"38\u00f75=7, 3\"],\n  [\"93\u00f79=10, 3\", \"21\u00f77=3, 0\"],\n  [\"75\u00f77=10, 5\", \"65\u00f72=32, 1\"],\n  [\"97\u00f77=13, 6\", \"20\u00f79=2, 2\"],\n  [\"24\u00f78=3, 0\", \"23\u00f76=3, 5\"],\n  [\"83\u00f72=41, 1\", \"59\u00f76=9, 5\"],\n  [\"25\u00f73=8, 1\", \"90\u00f77=12, 6\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  results.items[0].insertText(newText, 'Replace');\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-08-25 Sunday\", \"2024-08-26 Monday\")\n    ,@(\"37\u00f76=6, 1\", \"31\u00f76=5, 1\")\n    ,@(\"67\u00f79=7, 4\", \"57\u00f72=28, 1\")\n    ,@(\"10\u00f75=2, 0\", \"43\u00f75=8, 3\")\n    ,@(\"60\u00f76=10, 0\", \"93\u00f77=13, 2\")\n    ,@(\"66\u00f72=33, 0\", \"14\u00f76=2, 2\")\n    ,@(\"30\u00f76=5, 0\", \"97\u00f73=32, 1\")\n    ,@(\"88\u00f76=14, 4\", \"60\u00f77=8, 4\")\n    ,@(\"69\u00f77=9, 6\", \"44\u00f72=22, 0\")\n    ,@(\"15\u00f79=1, 6\", \"39\u00f76=6, 3\")\n    ,@(\"75\u00f74=18, 3\", \"42\u00f76=7, 0\")\n    ,@(\"98\u00f73=32, 2\", \"32\u00f75=6, 2\")\n    ,@(\"24\u00f73=8, 0\", \"95\u00f74=23, 3\")\n    ,@(\"79\u00f73=26, 1\", \"33\u00f78=4, 1\")\n    ,@(\"81\u00f79=9, 0\", \"13\u00f72=6, 1\")\n    ,@(\"43\u00f72=21, 1\", \"12\u00f76=2, 0\")\n    ,@(\"99\u00f72=49, 1\", \"47\u00f72=23, 1\")\n    ,@(\"51\u00f74=12, 3\", \"71\u00f78=8, 7\")\n    ,@(\"77\u00f72=38, 1\", \"65\u00f79=7, 2\")\n    ,@(\"43\u00f77=6, 1\", \"38\u00f75=7, 3\")\n    ,@(\"93\u00f79=10, 3\", \"21\u00f77=3, 0\")\n    ,@(\"75\u00f77=10, 5\", \"65\u00f72=32, 1\")\n    ,@(\"97\u00f77=13, 6\", \"20\u00f79=2, 2\")\n    ,@(\"24\u00f78=3, 0\", \"23\u00f76=3, 5\")\n    ,@(\"83\u00f72=41, 1\", \"59\u00f76=9, 5\")\n    ,@(\"25\u00f73=8, 1\", \"90\u00f77=12, 6\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $found = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}"}
